# Auto-generated Excel COM-interop script
# Updates cached market-price / profit values in the Pandaemonium_Profits workbook
# across several crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 206.25
$ws.Range("I8").Value = 206.25
$ws.Range("K8").Value = 618.75
$ws.Range("M8").Value = -479.75

$ws.Range("H33").Value = 4542450
$ws.Range("I33").Value = 6245792
$ws.Range("J33").Value = 203.66667
$ws.Range("K33").Value = 6245792
$ws.Range("L33").Value = 203.66667
$ws.Range("M33").Value = -6245563
$ws.Range("N33").Value = -661.6666700000001

$ws.Range("H43").Value = 1235.25
$ws.Range("J43").Value = 1350.2858
$ws.Range("L43").Value = 1350.2858
$ws.Range("N43").Value = -1488.2858

$ws.Range("H98").Value = 1986.2858
$ws.Range("I98").Value = 1394.4706
$ws.Range("K98").Value = 1394.4706
$ws.Range("M98").Value = 103.5293999999999

$ws.Range("H101").Value = 1918.5
$ws.Range("J101").Value = 2997.5
$ws.Range("L101").Value = 8992.5
$ws.Range("N101").Value = -12236.5

$ws.Range("H116").Value = 2128.5
$ws.Range("I116").Value = 1812.6666
$ws.Range("J116").Value = 2559.182
$ws.Range("K116").Value = 1812.6666
$ws.Range("L116").Value = 2559.182
$ws.Range("M116").Value = 1629.3334
$ws.Range("N116").Value = -9443.182000000001

$ws.Range("H122").Value = 1986.2858
$ws.Range("I122").Value = 1394.4706
$ws.Range("K122").Value = 4183.4118
$ws.Range("M122").Value = -1733.4118

$ws.Range("H137").Value = 656293.8
$ws.Range("I137").Value = 3914.8
$ws.Range("J137").Value = 928118.4399999999
$ws.Range("K137").Value = 11744.4
$ws.Range("L137").Value = 2784355.32
$ws.Range("M137").Value = -9194.400000000001
$ws.Range("N137").Value = -2789455.32

$ws.Range("H138").Value = 3336.28
$ws.Range("I138").Value = 2641.6667
$ws.Range("J138").Value = 3488.756
$ws.Range("K138").Value = 7925.000100000001
$ws.Range("L138").Value = 10466.268
$ws.Range("M138").Value = -2785.000100000001
$ws.Range("N138").Value = -20746.268

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 997.59375
$ws.Range("I2").Value = 1020.4643
$ws.Range("J2").Value = 837.5
$ws.Range("K2").Value = 1020.4643
$ws.Range("L2").Value = 837.5
$ws.Range("M2").Value = -907.4643
$ws.Range("N2").Value = -1063.5

$ws.Range("H32").Value = 16293.328
$ws.Range("I32").Value = 18025.951
$ws.Range("J32").Value = 2865.5
$ws.Range("K32").Value = 18025.951
$ws.Range("L32").Value = 2865.5
$ws.Range("M32").Value = -17738.951
$ws.Range("N32").Value = -3439.5

$ws.Range("H116").Value = 997.59375
$ws.Range("I116").Value = 1020.4643
$ws.Range("J116").Value = 837.5
$ws.Range("K116").Value = 1020.4643
$ws.Range("L116").Value = 837.5
$ws.Range("M116").Value = 1273.5357
$ws.Range("N116").Value = -5425.5

$ws.Range("H122").Value = 7354778
$ws.Range("I122").Value = 1930.2858
$ws.Range("J122").Value = 41668068
$ws.Range("K122").Value = 5790.857400000001
$ws.Range("L122").Value = 125004204
$ws.Range("M122").Value = -3340.857400000001
$ws.Range("N122").Value = -125009104

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 997.59375
$ws.Range("I3").Value = 1020.4643
$ws.Range("J3").Value = 837.5
$ws.Range("K3").Value = 1020.4643
$ws.Range("L3").Value = 837.5
$ws.Range("M3").Value = -906.4643
$ws.Range("N3").Value = -1065.5

$ws.Range("H99").Value = 3292.7646
$ws.Range("I99").Value = 1397.4445
$ws.Range("J99").Value = 5425
$ws.Range("K99").Value = 1397.4445
$ws.Range("L99").Value = 5425
$ws.Range("M99").Value = 100.5554999999999
$ws.Range("N99").Value = -8421

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1978.3158
$ws.Range("I16").Value = 1208
$ws.Range("J16").Value = 2671.6
$ws.Range("K16").Value = 1208
$ws.Range("L16").Value = 2671.6
$ws.Range("M16").Value = -921
$ws.Range("N16").Value = -3245.6

$ws.Range("H31").Value = 511466.22
$ws.Range("I31").Value = 7878.6787
$ws.Range("J31").Value = 824809.6
$ws.Range("K31").Value = 7878.6787
$ws.Range("L31").Value = 824809.6
$ws.Range("M31").Value = -7583.6787
$ws.Range("N31").Value = -825399.6

$ws.Range("H34").Value = 511466.22
$ws.Range("I34").Value = 7878.6787
$ws.Range("J34").Value = 824809.6
$ws.Range("K34").Value = 7878.6787
$ws.Range("L34").Value = 824809.6
$ws.Range("M34").Value = -7676.6787
$ws.Range("N34").Value = -825213.6

$ws.Range("H63").Value = 40271
$ws.Range("J63").Value = 40271
$ws.Range("L63").Value = 40271
$ws.Range("N63").Value = -41643

$ws.Range("H66").Value = 40271
$ws.Range("J66").Value = 40271
$ws.Range("L66").Value = 120813
$ws.Range("N66").Value = -127677

$ws.Range("H113").Value = 1978.3158
$ws.Range("I113").Value = 1208
$ws.Range("J113").Value = 2671.6
$ws.Range("K113").Value = 1208
$ws.Range("L113").Value = 2671.6
$ws.Range("M113").Value = 962
$ws.Range("N113").Value = -7011.6

$ws.Range("H132").Value = 2661.3052
$ws.Range("I132").Value = 1942.0488
$ws.Range("J132").Value = 4299.6113
$ws.Range("K132").Value = 5826.1464
$ws.Range("L132").Value = 12898.8339
$ws.Range("M132").Value = -3296.1464
$ws.Range("N132").Value = -17958.8339

$ws.Range("H134").Value = 1717.0927
$ws.Range("I134").Value = 1389.275
$ws.Range("J134").Value = 2653.7144
$ws.Range("K134").Value = 4167.825000000001
$ws.Range("L134").Value = 7961.1432
$ws.Range("M134").Value = -1632.825000000001
$ws.Range("N134").Value = -13031.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6176924
$ws.Range("I5").Value = 383.5
$ws.Range("J5").Value = 23824184
$ws.Range("K5").Value = 1150.5
$ws.Range("L5").Value = 71472552
$ws.Range("M5").Value = -1038.5
$ws.Range("N5").Value = -71472776

$ws.Range("H68").Value = 3601.3794
$ws.Range("I68").Value = 1588.8948
$ws.Range("J68").Value = 7425.1
$ws.Range("K68").Value = 4766.6844
$ws.Range("L68").Value = 22275.3
$ws.Range("M68").Value = -3955.6844
$ws.Range("N68").Value = -23897.3

$ws.Range("H71").Value = 3601.3794
$ws.Range("I71").Value = 1588.8948
$ws.Range("J71").Value = 7425.1
$ws.Range("K71").Value = 14300.0532
$ws.Range("L71").Value = 66825.90000000001
$ws.Range("M71").Value = -10244.0532
$ws.Range("N71").Value = -74937.90000000001

$ws.Range("H126").Value = 2674.6
$ws.Range("I126").Value = 1131.2858
$ws.Range("J126").Value = 4025
$ws.Range("K126").Value = 3393.8574
$ws.Range("L126").Value = 12075
$ws.Range("M126").Value = 1546.1426
$ws.Range("N126").Value = -21955

$ws.Range("H132").Value = 3114.6428
$ws.Range("I132").Value = 2740
$ws.Range("J132").Value = 3322.7778
$ws.Range("K132").Value = 24660
$ws.Range("L132").Value = 29905.0002
$ws.Range("M132").Value = -22130
$ws.Range("N132").Value = -34965.00019999999

$ws.Range("H133").Value = 3419.8215
$ws.Range("I133").Value = 1763.2354
$ws.Range("J133").Value = 5980
$ws.Range("K133").Value = 5289.706200000001
$ws.Range("L133").Value = 17940
$ws.Range("M133").Value = -229.7062000000005
$ws.Range("N133").Value = -28060

$ws.Range("H135").Value = 6176924
$ws.Range("I135").Value = 383.5
$ws.Range("J135").Value = 23824184
$ws.Range("K135").Value = 3451.5
$ws.Range("L135").Value = 214417656
$ws.Range("M135").Value = -916.5
$ws.Range("N135").Value = -214422726

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2447.1482
$ws.Range("I126").Value = 1805.2354
$ws.Range("J126").Value = 3538.4
$ws.Range("K126").Value = 5415.706200000001
$ws.Range("L126").Value = 10615.2
$ws.Range("M126").Value = -2945.706200000001
$ws.Range("N126").Value = -15555.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2520.087
$ws.Range("I7").Value = 2629.7856
$ws.Range("J7").Value = 2349.4443
$ws.Range("K7").Value = 2629.7856
$ws.Range("L7").Value = 2349.4443
$ws.Range("M7").Value = -2517.7856
$ws.Range("N7").Value = -2573.4443

$ws.Range("H40").Value = 3485.0454
$ws.Range("I40").Value = 3278.5715
$ws.Range("J40").Value = 3846.375
$ws.Range("K40").Value = 3278.5715
$ws.Range("L40").Value = 3846.375
$ws.Range("M40").Value = -3142.5715
$ws.Range("N40").Value = -4118.375

$ws.Range("H126").Value = 2520.087
$ws.Range("I126").Value = 2629.7856
$ws.Range("J126").Value = 2349.4443
$ws.Range("K126").Value = 7889.3568
$ws.Range("L126").Value = 7048.3329
$ws.Range("M126").Value = -5419.3568
$ws.Range("N126").Value = -11988.3329
